# Updated cryptos list (Price / Volume(1h) columns) to match the
# latest coinranking.com snapshot; also fixes the BinanceUSD/TRON
# row order (rows 16-17) which swapped position in the source feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.494.46"
$ws.Range("E2").Value = "  +0.73%  "

$ws.Range("D3").Value = "1.921.28"
$ws.Range("E3").Value = "  +1.20%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.57%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.63"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.54%  "

$ws.Range("E6").Value = "  +0.49%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4843"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +2.82%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4088"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.49%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08175"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.14%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.026"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.07%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.78"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +5.41%  "

$ws.Range("D12").Value = "1.889.83"
$ws.Range("E12").Value = "  -0.96%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.047"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +3.18%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.241"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.73%  "

$ws.Range("E15").Value = "  +2.89%  "

$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.007"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.72%  "

$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.06753"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.05%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001039"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.17%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.77"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.48%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.005"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.41%  "

$ws.Range("D21").Value = "29.525.48"
$ws.Range("E21").Value = "  +0.87%  "

$ws.Range("E22").Value = "  +2.41%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.77"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.95%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.180"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.83%  "

$ws.Range("D25").Value = "2.161.79"
$ws.Range("E25").Value = "  +2.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.749"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +10.79%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.63"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.47%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.11"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.21%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.125"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.58%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.66"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.84%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.024"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.46%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09561"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.14%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.543"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.58%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.569"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.72%  "

$ws.Range("E35").Value = "  -0.26%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02289"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.90%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06146"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.99%  "

$ws.Range("E38").Value = "  +0.98%  "

$ws.Range("E39").Value = "  +2.93%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.86"
$ws.Range("D40").ClearFormats()

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.014"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.78%  "

$ws.Range("E42").Value = "  +1.88%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.437"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.74%  "

$ws.Range("E44").Value = "  +1.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.07623"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.03%  "

$ws.Range("E46").Value = "  +2.93%  "

$ws.Range("E47").Value = "  +2.02%  "

$ws.Range("E48").Value = "  +3.21%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "116.87"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.432"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +4.29%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.82"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.63%  "
